$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -7.125399999999992
$ws.Range("D4").Value = -6.851399999999999
$ws.Range("C7").Value = -13.18929999999999
$ws.Range("A8").Value = -22.24790000000001
$ws.Range("A10").Value = -21.51099999999999
$ws.Range("E10").Value = 15.97349999999999
$ws.Range("D11").Value = -7.452099999999999
$ws.Range("A12").Value = -21.51170000000002
$ws.Range("E12").Value = 17.9692
$ws.Range("E13").Value = 16.33670000000001
$ws.Range("C14").Value = -13.447
$ws.Range("D14").Value = -8.366700000000005
$ws.Range("E14").Value = 15.95560000000001
$ws.Range("C15").Value = -14.37949999999998
$ws.Range("A18").Value = -22.0245
$ws.Range("C18").Value = -13.08579999999999
$ws.Range("D18").Value = -7.9132
$ws.Range("D19").Value = -9.223499999999987
$ws.Range("C20").Value = -11.9996
$ws.Range("D21").Value = -8.534999999999993
$ws.Range("A25").Value = -21.60669999999998
$ws.Range("D27").Value = -8.627100000000004
$ws.Range("C29").Value = -11.48690000000001
$ws.Range("E29").Value = 17.21450000000002
$ws.Range("C30").Value = -12.2872
$ws.Range("C31").Value = -12.4915
$ws.Range("D31").Value = -8.569500000000009
$ws.Range("E32").Value = 15.9726
$ws.Range("C35").Value = -12.1037
$ws.Range("E35").Value = 16.74070000000001
$ws.Range("A37").Value = -20.76930000000002
$ws.Range("D38").Value = -8.518700000000001
$ws.Range("C40").Value = -13.9836
$ws.Range("D42").Value = -9.004399999999993
$ws.Range("E43").Value = 17.6707
$ws.Range("C44").Value = -13.25949999999999
$ws.Range("D44").Value = -7.817800000000001
$ws.Range("D47").Value = -7.6787
$ws.Range("E48").Value = 17.49910000000002
$ws.Range("E49").Value = 15.93149999999998
$ws.Range("C50").Value = -13.6999
$ws.Range("E50").Value = 16.4254
$ws.Range("E51").Value = 17.30810000000001
$ws.Range("C54").Value = -13.02809999999999
$ws.Range("A55").Value = -22.22079999999999
$ws.Range("D56").Value = -8.665000000000001
$ws.Range("E56").Value = 16.08390000000001
$ws.Range("D58").Value = -8.226699999999996
$ws.Range("E61").Value = 16.64440000000001
$ws.Range("D65").Value = -7.926800000000002
$ws.Range("A68").Value = -21.71309999999999
$ws.Range("C68").Value = -11.7955
$ws.Range("E69").Value = 17.46760000000003
$ws.Range("E71").Value = 16.4805
$ws.Range("D73").Value = -7.687099999999996
$ws.Range("C76").Value = -12.1128
$ws.Range("A77").Value = -20.76179999999998
$ws.Range("A78").Value = -20.9971
$ws.Range("A79").Value = -20.80419999999998
$ws.Range("E79").Value = 18.51070000000001
$ws.Range("A80").Value = -20.5682
$ws.Range("A81").Value = -21.8398
$ws.Range("E81").Value = 16.40169999999999
$ws.Range("A82").Value = -21.81789999999999
$ws.Range("A84").Value = -22.0049
$ws.Range("C87").Value = -13.65409999999999
$ws.Range("C88").Value = -13.10529999999999
$ws.Range("D90").Value = -8.1304
$ws.Range("C92").Value = -12.4685
$ws.Range("D92").Value = -8.109700000000002
$ws.Range("E92").Value = 16.35659999999999
$ws.Range("D94").Value = -6.580999999999998
$ws.Range("D95").Value = -7.694200000000002
$ws.Range("C96").Value = -12.88070000000001
$ws.Range("C98").Value = -11.74499999999999
$ws.Range("A101").Value = -21.09259999999998
$ws.Range("C101").Value = -13.5702
$ws.Range("D101").Value = -7.992700000000002
$ws.Range("A102").Value = -19.44579999999998
$ws.Range("C102").Value = -13.31190000000001
